$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $text)
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextCell "D2" "27.919.38"
Set-TextCell "E2" "  +1.50%  "
Set-TextCell "D3" "1.639.85"
Set-TextCell "E3" "  +1.08%  "
Set-TextCell "E4" "  +0.03%  "
Set-TextCell "D5" "213.74"
Set-TextCell "E5" "  +0.97%  "
Set-TextCell "E7" "  +0.01%  "
Set-TextCell "D8" "23.63"
Set-TextCell "E8" "  +1.94%  "
Set-TextCell "E9" "  -0.27%  "
Set-TextCell "E10" "  +0.76%  "
Set-TextCell "D11" "0.0874"
Set-TextCell "E11" "  -0.88%  "
Set-TextCell "D12" "1.872.94"
Set-TextCell "E12" "  +1.14%  "
Set-TextCell "D13" "1.636.00"
Set-TextCell "E13" "  +0.76%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D14" "0.576"
Set-TextCell "E14" "  +4.57%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D15" "4.10"
Set-TextCell "E15" "  +1.34%  "
Set-TextCell "D16" "66.15"
Set-TextCell "E16" "  +1.22%  "
Set-TextCell "D17" "27.906.44"
Set-TextCell "E17" "  +1.52%  "
Set-TextCell "D18" "232.38"
Set-TextCell "E18" "  +1.28%  "
Set-TextCell "D19" "0.0₃0723"
Set-TextCell "E19" "  +0.66%  "
Set-TextCell "E20" "  +0.51%  "
Set-TextCell "E21" "  -0.02%  "
Set-TextCell "D22" "10.85"
Set-TextCell "E22" "  +3.91%  "
Set-TextCell "E23" "  +0.65%  "
Set-TextCell "E24" "  -3.16%  "
Set-TextCell "D25" "151.89"
Set-TextCell "E25" "  +1.74%  "
Set-TextCell "D26" "6.91"
Set-TextCell "E26" "  +0.65%  "
Set-TextCell "D27" "15.73"
Set-TextCell "E27" "  +1.39%  "
Set-TextCell "E28" "  +0.33%  "
Set-TextCell "E29" "  +0.03%  "
Set-TextCell "E30" "  +1.26%  "
Set-TextCell "E31" "  +0.15%  "
Set-TextCell "E32" "  +2.11%  "
Set-TextCell "E33" "  +1.95%  "
Set-TextCell "D34" "1.416.57"
Set-TextCell "E34" "  -3.41%  "
Set-TextCell "E35" "  +1.63%  "
Set-TextCell "E36" "  +0.39%  "
Set-TextCell "D37" "0.892"
Set-TextCell "E37" "  +2.09%  "
Set-TextCell "E38" "  +0.27%  "
Set-TextCell "E39" "  +0.42%  "
Set-TextCell "D40" "0.914"
Set-TextCell "E40" "  -3.55%  "
Set-TextCell "E41" "  +0.92%  "
Set-TextCell "E42" "  -0.02%  "
Set-TextCell "D43" "66.73"
Set-TextCell "E43" "  -1.91%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D44" "1.83"
Set-TextCell "E44" "  +4.12%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D45" "5.43"
Set-TextCell "E45" "  +1.98%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D46" "2.20"
Set-TextCell "E46" "  +0.35%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell "D47" "1.781.53"
Set-TextCell "E47" "  +1.23%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D48" "88.25"
Set-TextCell "E48" "  +1.20%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D49" "0.100"
Set-TextCell "E49" "  +1.01%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D50" "0.0506"
Set-TextCell "E50" "  +0.64%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D51" "7.64"
Set-TextCell "E51" "  +0.30%  "
